$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G14").Value = 1.558460351833249
$ws.Range("H2:H14").Value = 0.9990000000000001
